$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row above the current row 316, shifting the
# historical rows (316-369) down to (317-370), and fill the new row
# with the latest weekly price data.
$ws.Rows("316:316").Insert()

$ws.Range("A316").Value = 8
$ws.Range("B316").Value = "Terminal La Palmera de La Serena"
$ws.Range("C316").Value = "Coquimbo"
$ws.Range("D316").Value = 44476
$ws.Range("E316").Value = 4
$ws.Range("F316").Value = 100112024
$ws.Range("G316").Value = "Choclo"
$ws.Range("H316").Value = "Dulce o Americano"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 540
$ws.Range("K316").Value = 41000
$ws.Range("L316").Value = 42000
$ws.Range("M316").Value = 41500
$ws.Range("N316").Value = "$/malla 70 unidades"
$ws.Range("O316").Value = "Región de Arica y Parinacota"
$ws.Range("P316").Value = 593
$ws.Range("Q316").Value = 70
$ws.Range("R316").Value = "Hortaliza"
